$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 30304922  # H62 (was 23811410)
$ws.Cells.Item(62, 9).Value = 41668210  # I62 (was 41668230)
$ws.Cells.Item(62, 10).Value = 2831.6667  # J62 (was 2315.8333)
$ws.Cells.Item(62, 11).Value = 41668210  # K62 (was 41668230)
$ws.Cells.Item(62, 12).Value = 2831.6667  # L62 (was 2315.8333)
$ws.Cells.Item(62, 13).Value = -41667586  # M62 (was -41667606)
$ws.Cells.Item(62, 14).Value = -4079.6667  # N62 (was -3563.8333)
$ws.Cells.Item(65, 8).Value = 30304922  # H65 (was 23811410)
$ws.Cells.Item(65, 9).Value = 41668210  # I65 (was 41668230)
$ws.Cells.Item(65, 10).Value = 2831.6667  # J65 (was 2315.8333)
$ws.Cells.Item(65, 11).Value = 208341050  # K65 (was 208341150)
$ws.Cells.Item(65, 12).Value = 14158.3335  # L65 (was 11579.1665)
$ws.Cells.Item(65, 13).Value = -208337930  # M65 (was -208338030)
$ws.Cells.Item(65, 14).Value = -20398.3335  # N65 (was -17819.1665)
$ws.Cells.Item(76, 8).Value = 4392038.5  # H76 (was 4172579)
$ws.Cells.Item(76, 9).Value = 6951573  # I76 (was 6417056)
$ws.Cells.Item(76, 11).Value = 6951573  # K76 (was 6417056)
$ws.Cells.Item(76, 13).Value = -6951258  # M76 (was -6416741)
$ws.Cells.Item(79, 8).Value = 4392038.5  # H79 (was 4172579)
$ws.Cells.Item(79, 9).Value = 6951573  # I79 (was 6417056)
$ws.Cells.Item(79, 11).Value = 6951573  # K79 (was 6417056)
$ws.Cells.Item(79, 13).Value = -6950481  # M79 (was -6415964)
$ws.Cells.Item(88, 8).Value = 7582.1113  # H88 (was 7851.647)
$ws.Cells.Item(88, 10).Value = 8718.532999999999  # J88 (was 9127)
$ws.Cells.Item(88, 12).Value = 8718.532999999999  # L88 (was 9127)
$ws.Cells.Item(88, 14).Value = -9530.532999999999  # N88 (was -9939)
$ws.Cells.Item(91, 8).Value = 7582.1113  # H91 (was 7851.647)
$ws.Cells.Item(91, 10).Value = 8718.532999999999  # J91 (was 9127)
$ws.Cells.Item(91, 12).Value = 8718.532999999999  # L91 (was 9127)
$ws.Cells.Item(91, 14).Value = -11526.533  # N91 (was -11935)
$ws.Cells.Item(92, 8).Value = 73100300  # H92 (was 86806110)
$ws.Cells.Item(92, 9).Value = 5556387.5  # I92 (was 3968891)
$ws.Cells.Item(92, 10).Value = 148149090  # J92 (was 666666700)
$ws.Cells.Item(92, 11).Value = 5556387.5  # K92 (was 3968891)
$ws.Cells.Item(92, 12).Value = 148149090  # L92 (was 666666700)
$ws.Cells.Item(92, 13).Value = -5555139.5  # M92 (was -3967643)
$ws.Cells.Item(92, 14).Value = -148151586  # N92 (was -666669196)

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 630.1111  # H2 (was 637.03845)
$ws.Cells.Item(2, 9).Value = 480  # I2 (was 481.25)
$ws.Cells.Item(2, 11).Value = 480  # K2 (was 481.25)
$ws.Cells.Item(2, 13).Value = -367  # M2 (was -368.25)
$ws.Cells.Item(45, 8).Value = 6298.4736  # H45 (was 7303.6875)
$ws.Cells.Item(45, 9).Value = 6587.278  # I45 (was 7303.6875)
$ws.Cells.Item(45, 10).Value = 1100  # J45 (was 0)
$ws.Cells.Item(45, 11).Value = 6587.278  # K45 (was 7303.6875)
$ws.Cells.Item(45, 12).Value = 1100  # L45 (was 0)
$ws.Cells.Item(45, 13).Value = -6210.278  # M45 (was -6926.6875)
$ws.Cells.Item(45, 14).Value = -1854  # N45 (was None)
$ws.Cells.Item(97, 8).Value = 1450.8334  # H97 (was 1503.2941)
$ws.Cells.Item(97, 9).Value = 1370.2858  # I97 (was 1432.6923)
$ws.Cells.Item(97, 11).Value = 1370.2858  # K97 (was 1432.6923)
$ws.Cells.Item(97, 13).Value = -874.2858000000001  # M97 (was -936.6922999999999)
$ws.Cells.Item(110, 8).Value = 1654.9524  # H110 (was 1388.871)
$ws.Cells.Item(110, 9).Value = 1762.5  # I110 (was 1264.8235)
$ws.Cells.Item(110, 10).Value = 1588.7693  # J110 (was 1539.5)
$ws.Cells.Item(110, 11).Value = 1762.5  # K110 (was 1264.8235)
$ws.Cells.Item(110, 12).Value = 1588.7693  # L110 (was 1539.5)
$ws.Cells.Item(110, 13).Value = 282.5  # M110 (was 780.1765)
$ws.Cells.Item(110, 14).Value = -5678.7693  # N110 (was -5629.5)
$ws.Cells.Item(116, 8).Value = 630.1111  # H116 (was 637.03845)
$ws.Cells.Item(116, 9).Value = 480  # I116 (was 481.25)
$ws.Cells.Item(116, 11).Value = 480  # K116 (was 481.25)
$ws.Cells.Item(116, 13).Value = 1814  # M116 (was 1812.75)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 630.1111  # H3 (was 637.03845)
$ws.Cells.Item(3, 9).Value = 480  # I3 (was 481.25)
$ws.Cells.Item(3, 11).Value = 480  # K3 (was 481.25)
$ws.Cells.Item(3, 13).Value = -366  # M3 (was -367.25)
$ws.Cells.Item(33, 8).Value = 5512  # H33 (was 3704.6)
$ws.Cells.Item(33, 9).Value = 0  # I33 (was 1500)
$ws.Cells.Item(33, 10).Value = 5512  # J33 (was 4255.75)
$ws.Cells.Item(33, 11).Value = 0  # K33 (was 1500)
$ws.Cells.Item(33, 12).Value = 5512  # L33 (was 4255.75)
$ws.Cells.Item(33, 13).ClearContents()  # M33 (was -1164)
$ws.Cells.Item(33, 14).Value = -6184  # N33 (was -4927.75)
$ws.Cells.Item(94, 8).Value = 1879.5714  # H94 (was 1927.2858)
$ws.Cells.Item(94, 9).Value = 1191.2  # I94 (was 1329.4)
$ws.Cells.Item(94, 10).Value = 2262  # J94 (was 2259.4443)
$ws.Cells.Item(94, 11).Value = 1191.2  # K94 (was 1329.4)
$ws.Cells.Item(94, 12).Value = 2262  # L94 (was 2259.4443)
$ws.Cells.Item(94, 13).Value = -740.2  # M94 (was -878.4000000000001)
$ws.Cells.Item(94, 14).Value = -3164  # N94 (was -3161.4443)
$ws.Cells.Item(99, 8).Value = 47620252  # H99 (was 50001244)
$ws.Cells.Item(99, 9).Value = 83334270  # I99 (was 90910040)
$ws.Cells.Item(99, 10).Value = 1560.8889  # J99 (was 1609.7778)
$ws.Cells.Item(99, 11).Value = 83334270  # K99 (was 90910040)
$ws.Cells.Item(99, 12).Value = 1560.8889  # L99 (was 1609.7778)
$ws.Cells.Item(99, 13).Value = -83332772  # M99 (was -90908542)
$ws.Cells.Item(99, 14).Value = -4556.8889  # N99 (was -4605.7778)
$ws.Cells.Item(105, 8).Value = 9264.414000000001  # H105 (was 11918)
$ws.Cells.Item(105, 9).Value = 13044.111  # I105 (was 18302.23)
$ws.Cells.Item(105, 10).Value = 3079.4546  # J105 (was 3618.5)
$ws.Cells.Item(105, 11).Value = 13044.111  # K105 (was 18302.23)
$ws.Cells.Item(105, 12).Value = 3079.4546  # L105 (was 3618.5)
$ws.Cells.Item(105, 13).Value = -11297.111  # M105 (was -16555.23)
$ws.Cells.Item(105, 14).Value = -6573.4546  # N105 (was -7112.5)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 901.17645  # H22 (was 1225)
$ws.Cells.Item(22, 9).Value = 802.2222  # I22 (was 1200)
$ws.Cells.Item(22, 10).Value = 1012.5  # J22 (was 1266.6666)
$ws.Cells.Item(22, 11).Value = 802.2222  # K22 (was 1200)
$ws.Cells.Item(22, 12).Value = 1012.5  # L22 (was 1266.6666)
$ws.Cells.Item(22, 13).Value = -452.2222  # M22 (was -850)
$ws.Cells.Item(22, 14).Value = -1712.5  # N22 (was -1966.6666)
$ws.Cells.Item(31, 8).Value = 3094.62  # H31 (was 2856.5576)
$ws.Cells.Item(31, 9).Value = 1673.5135  # I31 (was 1562.619)
$ws.Cells.Item(31, 10).Value = 7139.3076  # J31 (was 8291.1)
$ws.Cells.Item(31, 11).Value = 1673.5135  # K31 (was 1562.619)
$ws.Cells.Item(31, 12).Value = 7139.3076  # L31 (was 8291.1)
$ws.Cells.Item(31, 13).Value = -1378.5135  # M31 (was -1267.619)
$ws.Cells.Item(31, 14).Value = -7729.3076  # N31 (was -8881.1)
$ws.Cells.Item(34, 8).Value = 3094.62  # H34 (was 2856.5576)
$ws.Cells.Item(34, 9).Value = 1673.5135  # I34 (was 1562.619)
$ws.Cells.Item(34, 10).Value = 7139.3076  # J34 (was 8291.1)
$ws.Cells.Item(34, 11).Value = 1673.5135  # K34 (was 1562.619)
$ws.Cells.Item(34, 12).Value = 7139.3076  # L34 (was 8291.1)
$ws.Cells.Item(34, 13).Value = -1471.5135  # M34 (was -1360.619)
$ws.Cells.Item(34, 14).Value = -7543.3076  # N34 (was -8695.1)
$ws.Cells.Item(58, 8).Value = 1325.0333  # H58 (was 1388.9434)
$ws.Cells.Item(58, 9).Value = 1092.8948  # I58 (was 1148)
$ws.Cells.Item(58, 10).Value = 1726  # J58 (was 1857.4445)
$ws.Cells.Item(58, 11).Value = 1092.8948  # K58 (was 1148)
$ws.Cells.Item(58, 12).Value = 1726  # L58 (was 1857.4445)
$ws.Cells.Item(58, 13).Value = -889.8948  # M58 (was -945)
$ws.Cells.Item(58, 14).Value = -2132  # N58 (was -2263.4445)
$ws.Cells.Item(94, 8).Value = 4825.5186  # H94 (was 4473.5)
$ws.Cells.Item(94, 9).Value = 4822  # I94 (was 4529.5)
$ws.Cells.Item(94, 10).Value = 4827.9375  # J94 (was 4436.1665)
$ws.Cells.Item(94, 11).Value = 4822  # K94 (was 4529.5)
$ws.Cells.Item(94, 12).Value = 4827.9375  # L94 (was 4436.1665)
$ws.Cells.Item(94, 13).Value = -4371  # M94 (was -4078.5)
$ws.Cells.Item(94, 14).Value = -5729.9375  # N94 (was -5338.1665)
$ws.Cells.Item(132, 8).Value = 2155.6978  # H132 (was 2581.5454)
$ws.Cells.Item(132, 9).Value = 1413.125  # I132 (was 1886.5)
$ws.Cells.Item(132, 11).Value = 4239.375  # K132 (was 5659.5)
$ws.Cells.Item(132, 13).Value = -1709.375  # M132 (was -3129.5)
$ws.Cells.Item(136, 8).Value = 1325.0333  # H136 (was 1388.9434)
$ws.Cells.Item(136, 9).Value = 1092.8948  # I136 (was 1148)
$ws.Cells.Item(136, 10).Value = 1726  # J136 (was 1857.4445)
$ws.Cells.Item(136, 11).Value = 3278.6844  # K136 (was 3444)
$ws.Cells.Item(136, 12).Value = 5178  # L136 (was 5572.333500000001)
$ws.Cells.Item(136, 13).Value = -728.6844000000001  # M136 (was -894)
$ws.Cells.Item(136, 14).Value = -10278  # N136 (was -10672.3335)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 3846503  # H38 (was 4348108)
$ws.Cells.Item(38, 10).Value = 428.77777  # J38 (was 209.83333)
$ws.Cells.Item(38, 12).Value = 1286.33331  # L38 (was 629.49999)
$ws.Cells.Item(38, 14).Value = -1980.33331  # N38 (was -1323.49999)
$ws.Cells.Item(113, 8).Value = 926526.4  # H113 (was 833920.4)
$ws.Cells.Item(113, 9).Value = 1282654.8  # I113 (was 1136951.2)
$ws.Cells.Item(113, 10).Value = 592.8  # J113 (was 585.375)
$ws.Cells.Item(113, 11).Value = 3847964.4  # K113 (was 3410853.6)
$ws.Cells.Item(113, 12).Value = 1778.4  # L113 (was 1756.125)
$ws.Cells.Item(113, 13).Value = -3845794.4  # M113 (was -3408683.6)
$ws.Cells.Item(113, 14).Value = -6118.4  # N113 (was -6096.125)
$ws.Cells.Item(137, 8).Value = 9407.516  # H137 (was 8799.416999999999)
$ws.Cells.Item(137, 9).Value = 7113.1055  # I137 (was 6711.905)
$ws.Cells.Item(137, 10).Value = 12521.357  # J137 (was 11721.934)
$ws.Cells.Item(137, 11).Value = 21339.3165  # K137 (was 20135.715)
$ws.Cells.Item(137, 12).Value = 37564.071  # L137 (was 35165.802)
$ws.Cells.Item(137, 13).Value = -16239.3165  # M137 (was -15035.715)
$ws.Cells.Item(137, 14).Value = -47764.071  # N137 (was -45365.802)

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5333.6514  # H70 (was 5578.9)
$ws.Cells.Item(70, 9).Value = 5434.44  # I70 (was 5690.561)
$ws.Cells.Item(70, 10).Value = 5018.6875  # J70 (was 5070.222)
$ws.Cells.Item(70, 11).Value = 5434.44  # K70 (was 5690.561)
$ws.Cells.Item(70, 12).Value = 5018.6875  # L70 (was 5070.222)
$ws.Cells.Item(70, 13).Value = -5164.44  # M70 (was -5420.561)
$ws.Cells.Item(70, 14).Value = -5558.6875  # N70 (was -5610.222)
$ws.Cells.Item(73, 8).Value = 5333.6514  # H73 (was 5578.9)
$ws.Cells.Item(73, 9).Value = 5434.44  # I73 (was 5690.561)
$ws.Cells.Item(73, 10).Value = 5018.6875  # J73 (was 5070.222)
$ws.Cells.Item(73, 11).Value = 5434.44  # K73 (was 5690.561)
$ws.Cells.Item(73, 12).Value = 5018.6875  # L73 (was 5070.222)
$ws.Cells.Item(73, 13).Value = -4498.44  # M73 (was -4754.561)
$ws.Cells.Item(73, 14).Value = -6890.6875  # N73 (was -6942.222)
$ws.Cells.Item(80, 8).Value = 2498.75  # H80 (was 2469.4119)
$ws.Cells.Item(80, 9).Value = 2498.3333  # I80 (was 2456.6667)
$ws.Cells.Item(80, 11).Value = 2498.3333  # K80 (was 2456.6667)
$ws.Cells.Item(80, 13).Value = -1500.3333  # M80 (was -1458.6667)
$ws.Cells.Item(83, 8).Value = 2498.75  # H83 (was 2469.4119)
$ws.Cells.Item(83, 9).Value = 2498.3333  # I83 (was 2456.6667)
$ws.Cells.Item(83, 11).Value = 12491.6665  # K83 (was 12283.3335)
$ws.Cells.Item(83, 13).Value = -7499.666499999999  # M83 (was -7291.333500000001)
$ws.Cells.Item(97, 8).Value = 1116.9524  # H97 (was 1159.174)
$ws.Cells.Item(97, 9).Value = 1097.3684  # I97 (was 1102.381)
$ws.Cells.Item(97, 10).Value = 1303  # J97 (was 1755.5)
$ws.Cells.Item(97, 11).Value = 1097.3684  # K97 (was 1102.381)
$ws.Cells.Item(97, 12).Value = 1303  # L97 (was 1755.5)
$ws.Cells.Item(97, 13).Value = -601.3684000000001  # M97 (was -606.3810000000001)
$ws.Cells.Item(97, 14).Value = -2295  # N97 (was -2747.5)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 49634.145  # H7 (was 2021.8889)
$ws.Cells.Item(7, 9).Value = 73293.71000000001  # I7 (was 2024.625)
$ws.Cells.Item(7, 10).Value = 2315  # J7 (was 2000)
$ws.Cells.Item(7, 11).Value = 73293.71000000001  # K7 (was 2024.625)
$ws.Cells.Item(7, 12).Value = 2315  # L7 (was 2000)
$ws.Cells.Item(7, 13).Value = -73181.71000000001  # M7 (was -1912.625)
$ws.Cells.Item(7, 14).Value = -2539  # N7 (was -2224)
$ws.Cells.Item(46, 8).Value = 917.4  # H46 (was 981.3684)
$ws.Cells.Item(46, 9).Value = 749.8333  # I46 (was 979.4)
$ws.Cells.Item(46, 10).Value = 989.2143  # J46 (was 982.0714)
$ws.Cells.Item(46, 11).Value = 749.8333  # K46 (was 979.4)
$ws.Cells.Item(46, 12).Value = 989.2143  # L46 (was 982.0714)
$ws.Cells.Item(46, 13).Value = -561.8333  # M46 (was -791.4)
$ws.Cells.Item(46, 14).Value = -1365.2143  # N46 (was -1358.0714)
$ws.Cells.Item(55, 8).Value = 281.88  # H55 (was 350.2353)
$ws.Cells.Item(55, 9).Value = 252.82353  # I55 (was 318.63635)
$ws.Cells.Item(55, 10).Value = 343.625  # J55 (was 408.16666)
$ws.Cells.Item(55, 11).Value = 252.82353  # K55 (was 318.63635)
$ws.Cells.Item(55, 12).Value = 343.625  # L55 (was 408.16666)
$ws.Cells.Item(55, 13).Value = -79.82353000000001  # M55 (was -145.63635)
$ws.Cells.Item(55, 14).Value = -689.625  # N55 (was -754.16666)
$ws.Cells.Item(100, 8).Value = 1411.4445  # H100 (was 1412.7778)
$ws.Cells.Item(100, 9).Value = 1337.875  # I100 (was 1339.375)
$ws.Cells.Item(100, 11).Value = 1337.875  # K100 (was 1339.375)
$ws.Cells.Item(100, 13).Value = -796.875  # M100 (was -798.375)
$ws.Cells.Item(126, 8).Value = 49634.145  # H126 (was 2021.8889)
$ws.Cells.Item(126, 9).Value = 73293.71000000001  # I126 (was 2024.625)
$ws.Cells.Item(126, 10).Value = 2315  # J126 (was 2000)
$ws.Cells.Item(126, 11).Value = 219881.13  # K126 (was 6073.875)
$ws.Cells.Item(126, 12).Value = 6945  # L126 (was 6000)
$ws.Cells.Item(126, 13).Value = -217411.13  # M126 (was -3603.875)
$ws.Cells.Item(126, 14).Value = -11885  # N126 (was -10940)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1320.2174  # H96 (was 1567.3684)
$ws.Cells.Item(96, 9).Value = 1437.7333  # I96 (was 1631.75)
$ws.Cells.Item(96, 10).Value = 1099.875  # J96 (was 1457)
$ws.Cells.Item(96, 11).Value = 1437.7333  # K96 (was 1631.75)
$ws.Cells.Item(96, 12).Value = 1099.875  # L96 (was 1457)
$ws.Cells.Item(96, 13).Value = -64.7333000000001  # M96 (was -258.75)
$ws.Cells.Item(96, 14).Value = -3845.875  # N96 (was -4203)
$ws.Cells.Item(136, 8).Value = 6469.2383  # H136 (was 5189.75)
$ws.Cells.Item(136, 9).Value = 45002  # I136 (was 11129.223)
$ws.Cells.Item(136, 10).Value = 2413.158  # J136 (was 2376.3157)
$ws.Cells.Item(136, 11).Value = 135006  # K136 (was 33387.669)
$ws.Cells.Item(136, 12).Value = 7239.474  # L136 (was 7128.9471)
$ws.Cells.Item(136, 13).Value = -132456  # M136 (was -30837.669)
$ws.Cells.Item(136, 14).Value = -12339.474  # N136 (was -12228.9471)
